# Updates the three-digit x one-digit multiplication problems in the
# document's tables. Each pair below is applied in an order chosen so
# that no intermediate replacement text collides with an old value that
# has not yet been processed (477x5= is both a source and a result).

$d = $word.ActiveDocument

$pairs = @(
    @("112×4=", "696×2="),
    @("193×4=", "999×8="),
    @("751×7=", "367×5="),
    @("227×7=", "272×3="),
    @("406×7=", "452×6="),
    @("908×9=", "414×2="),
    @("846×9=", "878×6="),
    @("982×7=", "714×2="),
    @("424×2=", "979×6="),
    @("477×5=", "537×7="),
    @("869×4=", "982×3="),
    @("234×2=", "246×8="),
    @("570×7=", "744×9="),
    @("378×7=", "215×2="),
    @("962×7=", "727×9="),
    @("650×2=", "248×7="),
    @("887×2=", "998×7="),
    @("160×8=", "193×7="),
    @("263×5=", "477×5="),
    @("564×2=", "621×8="),
    @("627×3=", "864×3="),
    @("620×9=", "136×3="),
    @("332×6=", "781×7="),
    @("585×4=", "922×4="),
    @("799×5=", "737×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
